$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "dni_ciu" (column E) values; column F ("PORC_AVANCE") is recomputed as E / D * 100
$updates = @(
    @{ Row = 2;  E = 947 }
    @{ Row = 3;  E = 763 }
    @{ Row = 4;  E = 942 }
    @{ Row = 5;  E = 872 }
    @{ Row = 6;  E = 928 }
    @{ Row = 7;  E = 952 }
    @{ Row = 8;  E = 662 }
    @{ Row = 9;  E = 1209 }
    @{ Row = 10; E = 791 }
    @{ Row = 11; E = 1430 }
    @{ Row = 12; E = 607 }
    @{ Row = 13; E = 869 }
    @{ Row = 14; E = 1945 }
    @{ Row = 15; E = 1308 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $dValue = $ws.Cells.Item($r, 4).Value2
    $eValue = $u.E
    $ws.Cells.Item($r, 5).Value2 = $eValue
    $ws.Cells.Item($r, 6).Value2 = ($eValue / $dValue) * 100
}
